$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row for the neural network result
$ws.Range("A7").Value = "NN (3 layers)"
$ws.Range("F7").Value = 0.77647942304611195

# Update the selection to match the saved workbook view state
$ws.Range("H14").Select()
